$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the shared string text "upb306" -> "Upb306" (cell A7)
$ws.Range("A7").Value = "Upb306"

# 2. Header row (A1:B1) gains the same border + font that the data rows already use
$ws.Range("A1:B1").Borders.LineStyle = $ws.Range("A2").Borders.LineStyle
$ws.Range("A1:B1").Font.Name = $ws.Range("A2").Font.Name
$ws.Range("A1:B1").Font.Size = $ws.Range("A2").Font.Size

# 3. Data rows (A2:B14) get an explicit black font color (was theme color before)
$ws.Range("A2:B14").Font.Color = 0

# 4. Data row heights grow from 18.75 to 19.5 (header row height is untouched)
$ws.Range("A2:B14").RowHeight = 19.5

Write-Output "done"
